# Update workbook per upstream data refresh (gh-pages output at 456a3b4)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 418
$ws1.Range("F4").Value  = 1176
$ws1.Range("F7").Value  = 41
$ws1.Range("F8").Value  = 1081
$ws1.Range("F9").Value  = 526
$ws1.Range("F10").Value = 386
$ws1.Range("F11").Value = 437
$ws1.Range("F13").Value = 323
$ws1.Range("F15").Value = 50
$ws1.Range("F17").Value = 17
$ws1.Range("F18").Value = 584
$ws1.Range("F19").Value = 1477
$ws1.Range("F20").Value = 5764
$ws1.Range("F22").Value = 1623
$ws1.Range("F23").Value = 386
$ws1.Range("F24").Value = 72
$ws1.Range("F25").Value = 35
$ws1.Range("F26").Value = 5394
$ws1.Range("F27").Value = 5394
$ws1.Range("F29").Value = 87
$ws1.Range("F30").Value = 1555
$ws1.Range("F31").Value = 91
$ws1.Range("F33").Value = 68
$ws1.Range("F34").Value = 1066
$ws1.Range("F35").Value = 672
$ws1.Range("F36").Value = 116
$ws1.Range("F37").Value = 4
$ws1.Range("F38").Value = 70
$ws1.Range("F39").Value = 3816

# ---------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 57
$ws2.Range("F5").Value = 174
$ws2.Range("G5").Value = 128
$ws2.Range("C6").Value = "杭州·2024吉卜力动漫音乐原版歌手交响音乐会（取消）"
$ws2.Range("G6").Value = "不可售"
$ws2.Range("F8").Value = 232
$ws2.Range("F19").Value = 26
$ws2.Range("F20").Value = 54

# ---------------------------------------------------------------
# Sheet: 本地生活 (Local life)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9429
$ws3.Range("F4").Value = 2162
$ws3.Range("F5").Value = 190

# ---------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9429
$ws4.Range("F4").Value  = 2162
$ws4.Range("F6").Value  = 418
$ws4.Range("F7").Value  = 1176
$ws4.Range("F10").Value = 41
$ws4.Range("F11").Value = 1081
$ws4.Range("F12").Value = 386
$ws4.Range("F13").Value = 437
$ws4.Range("F14").Value = 323
$ws4.Range("F16").Value = 50

# Row 18 now holds the "《卡农》永恒经典名曲音乐会" event (was "红楼梦" event)
$ws4.Range("C18").Value = "杭州·《卡农》永恒经典名曲音乐会"
$ws4.Range("D18").Value = "曙光路31号 浙江音乐厅"
$ws4.Range("E18").Value = "2024.09.07 19:30-09.07 21:00"
$ws4.Range("F18").Value = 2
$ws4.Range("G18").Value = 100
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=85894"
$ws4.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202405/3jz9YpaW1716100738530.jpeg"

# Row 19 now holds the "红楼梦" event (was "2024吉卜力" event)
# B19 holds a date-formatted *text* label ("2024-09-07"); writing it straight
# through .Value lets Excel auto-coerce it to a real date, so force text
# entry via NumberFormat and restore the Normal style afterwards so no
# spurious formatting is left behind.
$ws4.Range("B19").NumberFormat = "@"
$ws4.Range("B19").Value = "2024-09-07"
$ws4.Range("B19").Style = "Normal"
$ws4.Range("C19").Value = "杭州·红楼梦·主题演绎国风音乐会《梦寻红楼》"
$ws4.Range("D19").Value = "望梅路与汀兰路交叉口向南100米 杭州临平大剧院（原余杭大剧院）"
$ws4.Range("E19").Value = "2024.09.07 15:00-09.07 16:30"
$ws4.Range("F19").Value = 17
$ws4.Range("G19").Value = 100
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=89257"
$ws4.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202407/tkm6AHo71720572975141.jpeg"

$ws4.Range("F20").Value = 1478
$ws4.Range("F21").Value = 5764
$ws4.Range("F23").Value = 1623
$ws4.Range("F26").Value = 386
$ws4.Range("F29").Value = 5394
$ws4.Range("F30").Value = 5394
$ws4.Range("F32").Value = 87
$ws4.Range("F33").Value = 1555
$ws4.Range("F34").Value = 92
$ws4.Range("F36").Value = 1066
$ws4.Range("F37").Value = 672
$ws4.Range("F38").Value = 116
$ws4.Range("F44").Value = 70
$ws4.Range("F46").Value = 3816
$ws4.Range("F47").Value = 26
$ws4.Range("F48").Value = 54
